$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.139.19"
$ws.Range("E2").Value = "  +0.19%  "

$ws.Range("D3").Value = "1.676.15"
$ws.Range("E3").Value = "  -0.36%  "

$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.76%  "

$ws.Range("E6").Value = "  +0.15%  "

$ws.Range("E7").Value = "  +0.15%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.86"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +7.12%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.261"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.87%  "

$ws.Range("E10").Value = "  -0.47%  "

$ws.Range("E11").Value = "  +0.13%  "

$ws.Range("D12").Value = "1.914.22"
$ws.Range("E12").Value = "  -0.21%  "

$ws.Range("D13").Value = "1.680.98"
$ws.Range("E13").Value = "  +1.42%  "

$ws.Range("E14").Value = "  +2.13%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.561"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.54%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.52"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.10%  "

$ws.Range("D17").Value = "27.127.35"
$ws.Range("E17").Value = "  +0.22%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "235.17"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.51%  "

$ws.Range("E19").Value = "  +0.65%  "

$ws.Range("E20").Value = "  -4.09%  "

$ws.Range("E21").Value = "  +0.08%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.54"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.69%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.54"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.97%  "

$ws.Range("E24").Value = "  -1.99%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.81"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.26%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.47"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.37%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.45"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.68%  "

$ws.Range("E28").Value = "  -0.08%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0499"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.47%  "

$ws.Range("E31").Value = "  -0.21%  "

$ws.Range("E32").Value = "  -0.03%  "

$ws.Range("D33").Value = "1.545.97"
$ws.Range("E33").Value = "  -0.10%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.23"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.63%  "

$ws.Range("E35").Value = "  -3.98%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.608"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.53%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.947"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.45%  "

$ws.Range("E38").Value = "  +0.02%  "

$ws.Range("E39").Value = "  -1.06%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.07"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.26%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "69.72"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.69%  "

$ws.Range("E42").Value = "  +4.40%  "

$ws.Range("E43").Value = "  +0.13%  "

$ws.Range("E44").Value = "  -0.27%  "

$ws.Range("D45").Value = "1.823.70"
$ws.Range("E45").Value = "  +0.17%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.782"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.02%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.64"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.31%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "89.69"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.01%  "

$ws.Range("D49").Value = "0.0₆0111"
$ws.Range("E49").Value = "  +2.09%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.24"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.29%  "

$ws.Range("E51").Value = "  +0.35%  "
